$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StudentData")

# Append a new student record row (row 5)
$ws.Range("A5").Value = 20115817
$ws.Range("B5").Value = "Eli Mulaa"
$ws.Range("C5").Value = "eli@byupathway.edu"
$ws.Range("D5").Value = "Mulaaya@2025"
